$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 (I0) and J1 (IF), matching the style of H1
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in the data for columns I (I0) and J (IF) for rows 2-24
$values = @(
  9,9
  7,7
  6,7
  8,8
  8,9
  8,9
  6,7
  5,6
  9,9
  6,7
  6,6
  5,6
  6,6
  9,9
  8,8
  6,6
  7,7
  8,8
  6,6
  6,6
  4,4
  4,4
  5,5
)

for ($i = 0; $i -lt $values.Count; $i += 2) {
  $row = 2 + ($i / 2)
  $ws.Cells.Item($row, 9).Value = $values[$i]
  $ws.Cells.Item($row, 10).Value = $values[$i + 1]
}
